# Updates the cryptos list (Price / Volume(1h) columns) with refreshed
# market data. Plain numeric-looking prices are written with a leading
# apostrophe (forces Excel to keep them as text, matching the sheet's
# existing inline-string cells) and then restyled back to "Normal" so no
# stray quote-prefix style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.276.33"
$ws.Range("E2").Value = "  +8.76%  "

$ws.Range("D3").Value = "3.164.03"
$ws.Range("E3").Value = "  +6.40%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'592.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.14%  "

$ws.Range("D6").Value = "'148.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "3.161.32"
$ws.Range("E8").Value = "  +6.52%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("D10").Value = "'0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +18.77%  "

$ws.Range("E11").Value = "  +10.67%  "

$ws.Range("D12").Value = "'0.474"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.64%  "

$ws.Range("D13").Value = "'0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.24%  "

$ws.Range("D14").Value = "'36.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.23%  "

$ws.Range("D15").Value = "'0.124"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").Value = "3.680.09"
$ws.Range("E16").Value = "  +6.24%  "

$ws.Range("D17").Value = "'7.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.57%  "

$ws.Range("D18").Value = "64.099.67"
$ws.Range("E18").Value = "  +8.58%  "

$ws.Range("D19").Value = "3.148.55"
$ws.Range("E19").Value = "  +6.04%  "

$ws.Range("D20").Value = "'478.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.14%  "

$ws.Range("D21").Value = "'14.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.68%  "

$ws.Range("D22").Value = "'0.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "

$ws.Range("D23").Value = "'7.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.07%  "

$ws.Range("D24").Value = "'13.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("D25").Value = "'82.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.75%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").Value = "'8.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.19%  "

$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").Value = "'2.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.94%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").Value = "'6.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.89%  "

$ws.Range("D32").Value = "'27.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.44%  "

$ws.Range("D33").Value = "'0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.57%  "

$ws.Range("D34").Value = "0.0₃0894"
$ws.Range("E34").Value = "  +17.37%  "

$ws.Range("E35").Value = "  +19.14%  "

$ws.Range("E36").Value = "  +8.26%  "

$ws.Range("D37").Value = "'3.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +25.23%  "

$ws.Range("D38").Value = "'6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.55%  "

$ws.Range("D39").Value = "'51.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.64%  "

$ws.Range("D40").Value = "'452.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.18%  "

$ws.Range("D41").Value = "'8.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("D42").Value = "2.961.11"
$ws.Range("E42").Value = "  +8.76%  "

$ws.Range("D43").Value = "'0.0375"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.86%  "

$ws.Range("D44").Value = "'0.285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.33%  "

$ws.Range("D45").Value = "'0.112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.14%  "

$ws.Range("D46").Value = "'2.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +13.52%  "

$ws.Range("D47").Value = "'35.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.81%  "

$ws.Range("D49").Value = "'123.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("D51").Value = "'25.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.98%  "
